# Gaby time sheet and attendance updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Attendance: "Work (misc.)" row (row 14) now has hours logged on
# Thursday (E14) and Sunday (H14). Downstream SUM() formulas for the
# row/day/week totals recalc automatically.
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 1

# --- The Tue-Sun date header cells (C5:H5) were carrying a stray duplicate
# "bold" font (identical to the one already used by B5, just missing an
# explicit charset) left over from earlier editing. Re-apply the exact same
# font attributes as B5 so the workbook collapses back onto the single
# shared bold font/style instead of the redundant duplicate.
$hdr = $ws.Range("C5:H5")
$hdr.Font.Name = "Calibri"
$hdr.Font.Size = 11
$hdr.Font.Bold = $true
$hdr.Font.Color = 0

# --- Selection moved to H14 (last cell Gaby edited) before saving.
$ws.Range("H14").Select() | Out-Null
